$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.554.89"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.384.10"
$ws.Range("E3").Value = "  +6.92%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.639"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.17%  "
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.18%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.106"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").Value = "2.742.79"
$ws.Range("E16").Value = "  +6.94%  "
$ws.Range("D17").Value = "2.460.89"
$ws.Range("E17").Value = "  +9.63%  "
$ws.Range("D18").Value = "43.580.82"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  +5.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.64%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0925"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("E35").Value = "  +4.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0374"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.476"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.71%  "
